# Generate Report for Handoff
# A new handoff was generated for b.md (a new file
# b.63290e5768f688058c7b37413b0a5c26c308f864.<lang>.xlf), so the
# "Latest Handoff File" / "Latest Handoff Datetime" / "Status" for the
# b.md row on the zh-cn and de-de sheets (and the summary on Overview)
# need to reflect this new handoff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" row
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-24 20:42:20"

# ---------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" row
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-24 20:42:16"

# Rebuild the hyperlinks on the zh-cn sheet so the "display" text for D3
# points at the new handoff file, while every other hyperlink keeps its
# original address/display text.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/45c2e7be50f6c7b115af25afda5a6c532b874b3a/e2e/a.md", $null, $null, "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4fd0c922429f15ebf7c308b546ab7ed307f2c645/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7ae171218e20cd386d756dd61de565180f7eb280/e2e/a.md", $null, $null, "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f4a6db13bb4011ac43f38d723650342834ef8037/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/45c2e7be50f6c7b115af25afda5a6c532b874b3a/e2e/b.md", $null, $null, "b.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4fd0c922429f15ebf7c308b546ab7ed307f2c645/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", $null, $null, "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7ae171218e20cd386d756dd61de565180f7eb280/e2e/a.md", $null, $null, "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f4a6db13bb4011ac43f38d723650342834ef8037/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" row
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-24 20:42:20"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/45c2e7be50f6c7b115af25afda5a6c532b874b3a/e2e/a.md", $null, $null, "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f983606bf0385f1892ead680a99baafc3ba145b4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e37a68d683465c810027436dd074bf8002194145/e2e/a.md", $null, $null, "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a5aa477636ef659b35ffabf29014098ffa3ebef3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/45c2e7be50f6c7b115af25afda5a6c532b874b3a/e2e/b.md", $null, $null, "b.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f983606bf0385f1892ead680a99baafc3ba145b4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", $null, $null, "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e37a68d683465c810027436dd074bf8002194145/e2e/a.md", $null, $null, "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a5aa477636ef659b35ffabf29014098ffa3ebef3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", $null, $null, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
